$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("Q2").Value = 1.91
$ws.Range("R2").Value = 1.99
$ws.Range("W2").Value = 1.36

# Row 3
$ws.Range("AB3").Value = 8.5
$ws.Range("AL3").Value = 12
$ws.Range("AP3").Value = 34
$ws.Range("G3").Value = 1.85
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 4.5
$ws.Range("J3").Value = 2.5
$ws.Range("K3").Value = 2.1
$ws.Range("L3").Value = 4.75
$ws.Range("Q3").Value = 2.06
$ws.Range("R3").Value = 1.84
$ws.Range("W3").Value = 1.4
$ws.Range("Y3").Value = 1.8
$ws.Range("Z3").Value = 1.95

# Row 4
$ws.Range("G4").Value = 3.4
$ws.Range("I4").Value = 2.25
$ws.Range("K4").Value = 2
$ws.Range("U4").Value = 4.33
$ws.Range("V4").Value = 1.22

# Row 5
$ws.Range("AA5").Value = 5
$ws.Range("AD5").Value = 19
$ws.Range("AM5").Value = 19
$ws.Range("G5").Value = 2.05
$ws.Range("H5").Value = 2.82
$ws.Range("I5").Value = 4.2
$ws.Range("J5").Value = 3
$ws.Range("L5").Value = 5
$ws.Range("M5").Value = 1.14
$ws.Range("N5").Value = 5.5
$ws.Range("S5").Value = 5
$ws.Range("Z5").Value = 1.54

# Row 6
$ws.Range("AA6").Value = 6
$ws.Range("AB6").Value = 11
$ws.Range("AD6").Value = 26
$ws.Range("AE6").Value = 29
$ws.Range("AL6").Value = 6.5
$ws.Range("AM6").Value = 13
$ws.Range("AN6").Value = 12
$ws.Range("AO6").Value = 34
$ws.Range("G6").Value = 2.55
$ws.Range("H6").Value = 2.85
$ws.Range("I6").Value = 3.1
$ws.Range("L6").Value = 4
$ws.Range("Z6").Value = 1.54

# Row 7
$ws.Range("AR7").Value = 1.93
$ws.Range("AS7").Value = 1.93
$ws.Range("W7").Value = 1.57

# Row 8
$ws.Range("AA8").Value = 6
$ws.Range("AC8").Value = 10
$ws.Range("AF8").Value = 41
$ws.Range("AG8").Value = 6.5
$ws.Range("AL8").Value = 8.5
$ws.Range("AR8").Value = 1.88
$ws.Range("AS8").Value = 1.98
$ws.Range("K8").Value = 1.95
$ws.Range("Q8").Value = 2.5
$ws.Range("R8").Value = 1.5
$ws.Range("U8").Value = 5
$ws.Range("V8").Value = 1.17
$ws.Range("Y8").Value = 2.1
$ws.Range("Z8").Value = 1.67

# Row 9
$ws.Range("AG9").Value = 4.5
$ws.Range("AH9").Value = 5.5
$ws.Range("AM9").Value = 15
$ws.Range("G9").Value = 2.55
$ws.Range("I9").Value = 3.5
$ws.Range("K9").Value = 1.73
$ws.Range("M9").Value = 1.18
$ws.Range("N9").Value = 4.5
$ws.Range("O9").Value = 1.8
$ws.Range("P9").Value = 1.91
$ws.Range("Q9").Value = 3.6
$ws.Range("R9").Value = 1.29
$ws.Range("U9").Value = 9
$ws.Range("V9").Value = 1.07
$ws.Range("W9").Value = 1.83
$ws.Range("X9").Value = 1.98

# Row 13
$ws.Range("M13").Value = 1.03
$ws.Range("O13").Value = 1.17
$ws.Range("U13").Value = 2.62
$ws.Range("V13").Value = 1.41

# Row 14
$ws.Range("AA14").Value = 5
$ws.Range("AD14").Value = 9
$ws.Range("AG14").Value = 8
$ws.Range("AL14").Value = 13
$ws.Range("AQ14").Value = 67
$ws.Range("G14").Value = 1.45
$ws.Range("H14").Value = 4.33
$ws.Range("I14").Value = 7
$ws.Range("K14").Value = 2.38
$ws.Range("L14").Value = 7.5
$ws.Range("M14").Value = 1.05
$ws.Range("O14").Value = 1.37
$ws.Range("P14").Value = 2.75
$ws.Range("Q14").Value = 2.2
$ws.Range("R14").Value = 1.65
$ws.Range("V14").Value = 1.19
$ws.Range("W14").Value = 1.36
$ws.Range("X14").Value = 3
$ws.Range("Y14").Value = 2.5
$ws.Range("Z14").Value = 1.5

# Row 15
$ws.Range("S15").Value = 2.95
$ws.Range("T15").Value = 1.38

# Row 17
$ws.Range("AA17").Value = 7
$ws.Range("AC17").Value = 9.5
$ws.Range("AG17").Value = 7.5
$ws.Range("AH17").Value = 5.5
$ws.Range("AI17").Value = 15
$ws.Range("AK17").Value = 301
$ws.Range("AL17").Value = 10
$ws.Range("AM17").Value = 19
$ws.Range("G17").Value = 2.1
$ws.Range("H17").Value = 2.88
$ws.Range("I17").Value = 3.7
$ws.Range("J17").Value = 2.88
$ws.Range("K17").Value = 2
$ws.Range("M17").Value = 1.08
$ws.Range("N17").Value = 7.5
$ws.Range("Q17").Value = 2.25
$ws.Range("R17").Value = 1.62

# Row 18
$ws.Range("AG18").Value = 10
$ws.Range("AQ18").Value = 29
$ws.Range("G18").Value = 2.88
$ws.Range("K18").Value = 2.1
$ws.Range("M18").Value = 1.06
$ws.Range("N18").Value = 10
$ws.Range("Q18").Value = 1.98
$ws.Range("R18").Value = 1.88
$ws.Range("U18").Value = 3.4
$ws.Range("V18").Value = 1.3

# Row 28
$ws.Range("AM28").Value = 10
$ws.Range("AR28").Value = 1.9
$ws.Range("AS28").Value = 1.9
$ws.Range("G28").Value = 3.4
$ws.Range("I28").Value = 2.25
$ws.Range("J28").Value = 4
$ws.Range("M28").Value = 1.11
$ws.Range("N28").Value = 6.5
$ws.Range("O28").Value = 1.5
$ws.Range("P28").Value = 2.5
$ws.Range("V28").Value = 1.17
$ws.Range("W28").Value = 1.57
$ws.Range("X28").Value = 2.25
$ws.Range("Z28").Value = 1.69

# Row 29
$ws.Range("AB29").Value = 13
$ws.Range("AD29").Value = 29
$ws.Range("AI29").Value = 17
$ws.Range("AJ29").Value = 67
$ws.Range("AK29").Value = 800
$ws.Range("AL29").Value = 7
$ws.Range("AM29").Value = 12
$ws.Range("AN29").Value = 11
$ws.Range("AR29").Value = 1.85
$ws.Range("AS29").Value = 1.95
$ws.Range("G29").Value = 2.75
$ws.Range("H29").Value = 2.63
$ws.Range("I29").Value = 2.63
$ws.Range("J29").Value = 3.6
$ws.Range("K29").Value = 1.91
$ws.Range("L29").Value = 3.5
$ws.Range("M29").Value = 1.11
$ws.Range("N29").Value = 6.5
$ws.Range("O29").Value = 1.44
$ws.Range("P29").Value = 2.63
$ws.Range("Q29").Value = 2.5
$ws.Range("R29").Value = 1.5
$ws.Range("U29").Value = 5
$ws.Range("V29").Value = 1.17
$ws.Range("W29").Value = 1.57
$ws.Range("X29").Value = 2.25
$ws.Range("Z29").Value = 1.69

# Row 34
$ws.Range("AE34").Value = 26
$ws.Range("AG34").Value = 7.5
$ws.Range("AL34").Value = 6.5
$ws.Range("AM34").Value = 11
$ws.Range("AN34").Value = 11
$ws.Range("AO34").Value = 26
$ws.Range("G34").Value = 2.6
$ws.Range("I34").Value = 2.4
$ws.Range("J34").Value = 3.6
$ws.Range("L34").Value = 3.5
$ws.Range("M34").Value = 1.08
$ws.Range("N34").Value = 8
$ws.Range("Y34").Value = 2.1
$ws.Range("Z34").Value = 1.67

# Row 36
$ws.Range("Q36").Value = 1.83
$ws.Range("R36").Value = 2.03

# Row 37
$ws.Range("AR37").Value = 1.43
$ws.Range("AS37").Value = 2.85
$ws.Range("M37").Value = 1.04
$ws.Range("N37").Value = 13
$ws.Range("Q37").Value = 1.88
$ws.Range("R37").Value = 1.98
$ws.Range("T37").Value = 1.53

# Row 40
$ws.Range("AE40").Value = 17
$ws.Range("AN40").Value = 13
$ws.Range("G40").Value = 1.91
$ws.Range("I40").Value = 4
$ws.Range("J40").Value = 2.63
$ws.Range("L40").Value = 4.33
$ws.Range("Q40").Value = 2.03
$ws.Range("R40").Value = 1.83
$ws.Range("U40").Value = 3.4
$ws.Range("V40").Value = 1.3
$ws.Range("W40").Value = 1.44
$ws.Range("X40").Value = 2.63

# Row 48
$ws.Range("AB48").Value = 8.5
$ws.Range("AD48").Value = 17
$ws.Range("AG48").Value = 7.5
$ws.Range("AL48").Value = 8.5
$ws.Range("AM48").Value = 17
$ws.Range("AR48").Value = 1.8
$ws.Range("AS48").Value = 2.05
$ws.Range("G48").Value = 2.05
$ws.Range("H48").Value = 3.25
$ws.Range("I48").Value = 3.6
$ws.Range("J48").Value = 2.88
$ws.Range("K48").Value = 2
$ws.Range("L48").Value = 4.5
$ws.Range("N48").Value = 7.5
$ws.Range("Q48").Value = 2.35
$ws.Range("R48").Value = 1.57
$ws.Range("W48").Value = 1.53
$ws.Range("X48").Value = 2.38
